# New crime data collected - weekly update for week covering 11/13/2023 - 11/19/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates ----
# Volume/Number line: "... Number 45" -> "... Number 46"
$ws.Range("A8").Value = "Volume 30   Number  46"

# Report covering week line
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# ---- Donor cells used to flip a cell between "numeric" and "text-shared-string"
# representations while keeping the exact same style index (s="14" for text
# cells that hold the literal strings "0" / "***.*", s="15"/"16" for plain
# numeric cells). Row 14 is untouched by this update, so its cells make
# reliable donors:
#   C14 -> text "0"   (style 14)
#   E14 -> text "***.*" (style 14)
#   F14 -> plain number (style 15)
#   K14 -> plain number (style 16)

# ===================== Row 14 (Murder) =====================
$ws.Range("N14").Value = -82.692307692307

# ===================== Row 15 (Rape) =====================
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("M15").Value = -20
$ws.Range("N15").Value = -69.230769230769

# ===================== Row 16 (Robbery) =====================
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 186
$ws.Range("J16").Value = 196
$ws.Range("K16").Value = -5.102040816326
$ws.Range("L16").Value = 20.779220779220
$ws.Range("M16").Value = -27.058823529411
$ws.Range("N16").Value = -77.149877149877

# ===================== Row 17 (Fel. Assault) =====================
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -77.777777777777
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 42
$ws.Range("H17").Value = -45.238095238095
$ws.Range("I17").Value = 406
$ws.Range("J17").Value = 462
$ws.Range("K17").Value = -12.121212121212
$ws.Range("L17").Value = -13.247863247863
$ws.Range("M17").Value = 42.957746478873
$ws.Range("N17").Value = -54.535274356103

# ===================== Row 18 (Burglary) =====================
$ws.Range("F14").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -53.333333333333
$ws.Range("I18").Value = 152
$ws.Range("J18").Value = 179
$ws.Range("K18").Value = -15.083798882681
$ws.Range("L18").Value = 0.662251655629
$ws.Range("M18").Value = 34.513274336283
$ws.Range("N18").Value = -82.568807339449

# ===================== Row 19 (Gr. Larceny) =====================
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 35.714285714285
$ws.Range("I19").Value = 395
$ws.Range("J19").Value = 307
$ws.Range("K19").Value = 28.664495114006
$ws.Range("L19").Value = 39.575971731448
$ws.Range("M19").Value = 45.220588235294
$ws.Range("N19").Value = 9.722222222222

# ===================== Row 20 (G.L.A.) =====================
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 116.666666666667
$ws.Range("I20").Value = 129
$ws.Range("J20").Value = 67
$ws.Range("K20").Value = 92.537313432835
$ws.Range("L20").Value = 148.076923076923
$ws.Range("M20").Value = 180.434782608696
$ws.Range("N20").Value = -42.920353982300

# ===================== Row 21 (TOTAL) =====================
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -26.923076923076
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -10.526315789473
$ws.Range("I21").Value = 1301
$ws.Range("J21").Value = 1242
$ws.Range("K21").Value = 4.750402576489
$ws.Range("L21").Value = 13.922942206655
$ws.Range("M21").Value = 29.195630585898
$ws.Range("N21").Value = -60.51593323217

# ===================== Row 22 (Transit) =====================
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("F14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = -26.315789473684

# ===================== Row 23 (Housing) =====================
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 27
$ws.Range("H23").Value = -44.444444444444
$ws.Range("I23").Value = 214
$ws.Range("J23").Value = 208
$ws.Range("K23").Value = 2.884615384615
$ws.Range("L23").Value = 8.080808080808
$ws.Range("M23").Value = 41.721854304635

# ===================== Row 24 (Petit Larceny) =====================
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = 23.076923076923
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 50
$ws.Range("H24").Value = 32
$ws.Range("I24").Value = 783
$ws.Range("J24").Value = 733
$ws.Range("K24").Value = 6.821282401091
$ws.Range("L24").Value = 15.486725663716
$ws.Range("M24").Value = 36.649214659685

# ===================== Row 25 (Misd. Assault) =====================
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = -18.333333333333
$ws.Range("I25").Value = 543
$ws.Range("J25").Value = 539
$ws.Range("K25").Value = 0.742115027829
$ws.Range("L25").Value = -1.451905626134
$ws.Range("M25").Value = -29.019607843137

# ===================== Row 26 (UCR Rape*) =====================
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 33.333333333333
$ws.Range("I26").Value = 40
$ws.Range("K26").Value = -4.761904761904
$ws.Range("L26").Value = 29.032258064516

# ===================== Row 27 (Other Sex Crimes) =====================
$ws.Range("F14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 48
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = -26.153846153846
$ws.Range("L27").Value = -25

# ===================== Row 28 (Shooting Vic.) =====================
$ws.Range("L28").Value = -51.785714285714
$ws.Range("N28").Value = -77.310924369747

# ===================== Row 29 (Shooting Inc.) =====================
$ws.Range("L29").Value = -46.666666666666
$ws.Range("N29").Value = -77.777777777777

# ===================== Row 30 (Hate Crimes) =====================
$ws.Range("C14").Copy($ws.Range("G30"))
$ws.Range("E14").Copy($ws.Range("H30"))

Write-Host "Weekly crime data update applied"
